# "Update countries & provincias Spain" — refresh the COVID figures in the
# "Pais" sheet and bump the "Datos actualizados" timestamp.
#
# The underlying data table (rows 4..216) is emitted sorted by total cases,
# so as country totals move, whole rows of figures shift to new ranks while
# column A (country name) is rewritten to match. Kirguistan's case count
# grew enough to leapfrog Libano/Crucero/Niger (rows 98-101), and Pakistan's
# row (30) got fresh totals too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer cell (A1): 05:52 -> 06:22
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 29 de Abril de 2020 a las 06:22"

# Row 30 - Pakistan: refreshed totals
$ws.Cells.Item(30, 2).Value = 14885   # Casos totales
$ws.Cells.Item(30, 3).Value = 273     # Nuevos casos
$ws.Cells.Item(30, 4).Value = 3425    # Casos activos
$ws.Cells.Item(30, 5).Value = 11133   # Recuperados
$ws.Cells.Item(30, 6).Value = 111     # Casos criticos (unchanged)
$ws.Cells.Item(30, 7).Value = 15      # Muertes hoy
$ws.Cells.Item(30, 8).Value = 327     # Muertes

# Row 98 - Kirguistan moves up into this rank, with fresh totals
$ws.Cells.Item(98, 1).Value = "Kirguistan"
$ws.Cells.Item(98, 2).Value = 729
$ws.Cells.Item(98, 3).Value = 21
$ws.Cells.Item(98, 4).Value = 437
$ws.Cells.Item(98, 5).Value = 284
$ws.Cells.Item(98, 6).Value = 13
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 8

# Row 99 - Libano drops one rank (keeps its previous day's totals)
$ws.Cells.Item(99, 1).Value = "Libano"
$ws.Cells.Item(99, 2).Value = 717
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 145
$ws.Cells.Item(99, 5).Value = 548
$ws.Cells.Item(99, 6).Value = 44
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 24

# Row 100 - Crucero drops one rank
$ws.Cells.Item(100, 1).Value = "Crucero"
$ws.Cells.Item(100, 2).Value = 712
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 645
$ws.Cells.Item(100, 5).Value = 54
$ws.Cells.Item(100, 6).Value = 4
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 13

# Row 101 - Niger drops one rank, takes Kirguistan's old slot
$ws.Cells.Item(101, 1).Value = "Niger"
$ws.Cells.Item(101, 2).Value = 709
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 403
$ws.Cells.Item(101, 5).Value = 275
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 31
